$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "B2" = 0.2036036036036036
    "C2" = 0.5513513513513514
    "J2" = 0.01621621621621622
    "P2" = 0.1441441441441441
    "S2" = 0.08468468468468468
    "B3" = 0.003236245954692557
    "C3" = 0.02588996763754045
    "J3" = 0.03883495145631068
    "P3" = 0.7508090614886731
    "S3" = 0.1812297734627832
    "J4" = 0.05825242718446602
    "P4" = 0.7184466019417476
    "S4" = 0.2233009708737864
    "P5" = 0.3333333333333333
    "S5" = 0.6666666666666666
    "B6" = 0.036281179138322
    "D6" = 0.01133786848072562
    "E6" = 0.002267573696145125
    "F6" = 0.06122448979591837
    "J6" = 0.253968253968254
    "O6" = 0.018140589569161
    "Q6" = 0.1564625850340136
    "R6" = 0.08843537414965986
    "S6" = 0.3718820861678004
    "B7" = 0.09921671018276762
    "D7" = 0.02349869451697128
    "E7" = 0.005221932114882507
    "F7" = 0.04699738903394256
    "J7" = 0.1331592689295039
    "O7" = 0.01827676240208877
    "Q7" = 0.1566579634464752
    "R7" = 0.06788511749347259
    "S7" = 0.4490861618798956
    "B8" = 0.09582863585118377
    "D8" = 0.01578354002254791
    "E8" = 0.00112739571589628
    "F8" = 0.06313416009019165
    "J8" = 0.1172491544532131
    "O8" = 0.01240135287485907
    "Q8" = 0.1927846674182638
    "R8" = 0.1071025930101466
    "S8" = 0.3945885005636979
    "B9" = 0.07954545454545454
    "D9" = 0.02045454545454545
    "F9" = 0.06818181818181818
    "J9" = 0.1068181818181818
    "O9" = 0.01136363636363636
    "Q9" = 0.1863636363636364
    "R9" = 0.09545454545454546
    "S9" = 0.4318181818181818
    "B10" = 0.09438775510204081
    "D10" = 0.02551020408163265
    "E10" = 0.0003644314868804665
    "F10" = 0.06523323615160349
    "J10" = 0.1202623906705539
    "O10" = 0.01639941690962099
    "Q10" = 0.228134110787172
    "R10" = 0.08782798833819241
    "S10" = 0.3618804664723032
    "F11" = 0.001669449081803005
    "G11" = 0.1452420701168614
    "J11" = 0.09682804674457429
    "K11" = 0.1953255425709516
    "L11" = 0.5492487479131887
    "S11" = 0.01168614357262104
    "G12" = 0.7191011235955056
    "J12" = 0.1882022471910112
    "K12" = 0.01685393258426966
    "L12" = 0.06741573033707865
    "S12" = 0.008426966292134831
    "G13" = 0.6052631578947368
    "J13" = 0.3552631578947368
    "S13" = 0.03947368421052631
    "G14" = 0.625
    "J14" = 0.375
    "F15" = 0.0207852193995381
    "H15" = 0.1547344110854504
    "I15" = 0.07621247113163972
    "J15" = 0.3602771362586605
    "K15" = 0.07390300230946882
    "M15" = 0.002309468822170901
    "O15" = 0.06697459584295612
    "S15" = 0.2448036951501155
    "F16" = 0.01861702127659574
    "H16" = 0.1914893617021277
    "I16" = 0.09574468085106383
    "J16" = 0.3723404255319149
    "K16" = 0.1143617021276596
    "M16" = 0.02127659574468085
    "N16" = 0.002659574468085106
    "O16" = 0.05851063829787234
    "S16" = 0.125
    "F17" = 0.01592039800995025
    "H17" = 0.173134328358209
    "I17" = 0.09751243781094528
    "J17" = 0.4308457711442786
    "K17" = 0.09154228855721393
    "M17" = 0.01691542288557214
    "N17" = 0.0009950248756218905
    "O17" = 0.05970149253731343
    "S17" = 0.1134328358208955
    "F18" = 0.0137299771167048
    "H18" = 0.1510297482837529
    "I18" = 0.09610983981693363
    "J18" = 0.4622425629290618
    "K18" = 0.09382151029748284
    "M18" = 0.009153318077803204
    "N18" = 0.002288329519450801
    "O18" = 0.06178489702517163
    "S18" = 0.1098398169336384
    "F19" = 0.01439688715953307
    "H19" = 0.198443579766537
    "I19" = 0.08949416342412451
    "J19" = 0.3976653696498054
    "K19" = 0.1011673151750973
    "M19" = 0.01828793774319066
    "N19" = 0.002334630350194552
    "O19" = 0.0669260700389105
    "S19" = 0.111284046692607
}

foreach ($key in $values.Keys) {
    $ws.Range($key).Value = $values[$key]
}